$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
$c1 = $cs.Colors(1)
$c1.RGB = 42
